$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -6.798699999999992
$ws.Range("A3").Value = -21.20060000000003
$ws.Range("D5").Value = -8.453099999999996
$ws.Range("E5").Value = 12.71329999999999
$ws.Range("E9").Value = 13.53450000000001
$ws.Range("E11").Value = 13.4032
$ws.Range("A14").Value = -20.42079999999998
$ws.Range("A21").Value = -21.33500000000001
$ws.Range("E21").Value = 12.8393
$ws.Range("A23").Value = -21.39870000000003
$ws.Range("A25").Value = -22.40690000000003
